$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 106 - this shifts existing rows 106..205 down to 107..206,
# preserving all of their data/styles (matches the diff's observed "cascade" where
# every row from 106 downward takes on the values of the row above it).
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new record added by this edit.
$ws.Cells.Item(106, 1).Value2 = 10
$ws.Cells.Item(106, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value2 = "La Araucanía"
$ws.Cells.Item(106, 4).Value2 = 44484
$ws.Cells.Item(106, 5).Value2 = 9
$ws.Cells.Item(106, 6).Value2 = 100112044
$ws.Cells.Item(106, 7).Value2 = "Perejil"
$ws.Cells.Item(106, 8).Value2 = "Sin especificar"
$ws.Cells.Item(106, 9).Value2 = "Primera"
$ws.Cells.Item(106, 10).Value2 = 40
$ws.Cells.Item(106, 11).Value2 = 4000
$ws.Cells.Item(106, 12).Value2 = 4500
$ws.Cells.Item(106, 13).Value2 = 4250
$ws.Cells.Item(106, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(106, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(106, 16).Value2 = 1417
$ws.Cells.Item(106, 17).Value2 = 3
$ws.Cells.Item(106, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the rest of
# column D (same style as the cell above/below it).
$ws.Cells.Item(106, 4).NumberFormat = $ws.Cells.Item(105, 4).NumberFormat
